# adding profits to tables
# Adds a third block of 8 columns (R:Y) mirroring the existing M_%cit (B:I)
# and M_ETR (J:Q) blocks, labeled "M_PL", with new profit data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 1: new merged header "M_PL" over R1:Y1 (mirrors B1:I1 "M_%cit"
# and J1:Q1 "M_ETR")
# ---------------------------------------------------------------
$ws.Range("R1").Value = "M_PL"
$ws.Range("R1:Y1").Merge()

# ---------------------------------------------------------------
# Row 2: repeat the same 8 sub-headers used under the other two blocks
# ---------------------------------------------------------------
$headers = @("GFA - Sales", "GFA - Sales + Emp", "IMF - Sales", "IMF - Sales + Emp", "OECD (20%) - Sales", "OECD (20%) - Sales + Emp", "OECD - Sales", "OECD - Sales + Emp")
$cols = @("R", "S", "T", "U", "V", "W", "X", "Y")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $headers[$i]
}

# ---------------------------------------------------------------
# Match formatting of the existing header cells (bold, centered,
# top-aligned, thin border all round) for the new header cells.
# ---------------------------------------------------------------
$headerRange = $ws.Range("R1:Y2")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# ---------------------------------------------------------------
# Data rows 4-10: new profit figures for each of the 8 sub-columns
# ---------------------------------------------------------------
$ws.Range("R4").Value = 40382028217
$ws.Range("S4").Value = 37272812891
$ws.Range("T4").Value = 30234495459
$ws.Range("U4").Value = 37431215439
$ws.Range("V4").Value = 66677217402
$ws.Range("W4").Value = 66677217402
$ws.Range("X4").Value = 66677217402
$ws.Range("Y4").Value = 66677217402

$ws.Range("R5").Value = 59581902456
$ws.Range("S5").Value = 59847798241
$ws.Range("T5").Value = 57906933849
$ws.Range("U5").Value = 58157825465
$ws.Range("V5").Value = 66483805754
$ws.Range("W5").Value = 66483805754
$ws.Range("X5").Value = 66483805754
$ws.Range("Y5").Value = 66483805754

$ws.Range("R6").Value = 12716278503
$ws.Range("S6").Value = 12679471304
$ws.Range("T6").Value = 11918971016
$ws.Range("U6").Value = 12679471304
$ws.Range("V6").Value = 17853638274
$ws.Range("W6").Value = 17853638274
$ws.Range("X6").Value = 17853638274
$ws.Range("Y6").Value = 17853638274

$ws.Range("R7").Value = 699911008
$ws.Range("S7").Value = 699911008
$ws.Range("T7").Value = 699911008
$ws.Range("U7").Value = 699911008
$ws.Range("V7").Value = 5767699375
$ws.Range("W7").Value = 5767699375
$ws.Range("X7").Value = 5767699375
$ws.Range("Y7").Value = 5767699375

$ws.Range("R8").Value = 885447038872
$ws.Range("S8").Value = 885447038872
$ws.Range("T8").Value = 885447038872
$ws.Range("U8").Value = 885447038872
$ws.Range("V8").Value = 885447038872
$ws.Range("W8").Value = 885447038872
$ws.Range("X8").Value = 885447038872
$ws.Range("Y8").Value = 885447038872

# row 9 mirrors the sparse population pattern already present in B9:Q9
# (R9 and T9 are intentionally left blank)
$ws.Range("S9").Value = 12096909667
$ws.Range("U9").Value = 12096909667
$ws.Range("V9").Value = 12956669707
$ws.Range("W9").Value = 12956669707
$ws.Range("X9").Value = 12956669707
$ws.Range("Y9").Value = 12956669707

$ws.Range("R10").Value = 660204378
$ws.Range("S10").Value = -130919670
$ws.Range("T10").Value = 660204378
$ws.Range("U10").Value = 705481484
$ws.Range("V10").Value = 1212797515
$ws.Range("W10").Value = 1212797515
$ws.Range("X10").Value = 1212797515
$ws.Range("Y10").Value = 1212797515
